$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Wild thing" meeting note text for the relevant meeting (row 8)
$ws.Range("E8").Value = 'Wild thing finalised as "Taipan". Data pulled using the Galah API. Next steps discussed briefly.'

# Update the active selection to match the edited cell
$ws.Range("E8").Select()
